$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1113813209086354
$ws.Range("C2").Value = 1.082668180363185
$ws.Range("D2").Value = 3.265531177099384
$ws.Range("E2").Value = 1.807078077200701
$ws.Range("F2").Value = 1.828521252897521
$ws.Range("G2").Value = 37

$ws.Range("B3").Value = 0.1513751093172216
$ws.Range("C3").Value = 0.9545256880515617
$ws.Range("D3").Value = 2.919874515134196
$ws.Range("E3").Value = 1.7087640314374
$ws.Range("F3").Value = 1.726189554985052
$ws.Range("G3").Value = 36

$ws.Range("B4").Value = 0.2122263332368638
$ws.Range("C4").Value = 0.9540284968671798
$ws.Range("D4").Value = 2.851160657205658
$ws.Range("E4").Value = 1.688538023618555
$ws.Range("F4").Value = 1.699603954074797
$ws.Range("G4").Value = 35

$ws.Range("B5").Value = 0.1511539997107805
$ws.Range("C5").Value = 0.9530044890603636
$ws.Range("D5").Value = 2.951427543481194
$ws.Range("E5").Value = 1.717971927442702
$ws.Range("F5").Value = 1.737044864330408
$ws.Range("G5").Value = 34

$ws.Range("B6").Value = 0.2330575124652057
$ws.Range("C6").Value = 0.9893939392548374
$ws.Range("D6").Value = 2.998246469334817
$ws.Range("E6").Value = 1.731544532876593
$ws.Range("F6").Value = 1.742391603086522
$ws.Range("G6").Value = 33

$ws.Range("B7").Value = 0.1708011111517192
$ws.Range("C7").Value = 0.9843372303369592
$ws.Range("D7").Value = 3.106494255942373
$ws.Range("E7").Value = 1.762524966047963
$ws.Range("F7").Value = 1.782298982592832
$ws.Range("G7").Value = 32

$ws.Range("B8").Value = 0.2456821311865818
$ws.Range("C8").Value = 1.039250576912814
$ws.Range("D8").Value = 3.172757163278162
$ws.Range("E8").Value = 1.781223501775721
$ws.Range("F8").Value = 1.793361100136716
$ws.Range("G8").Value = 31

$ws.Range("B9").Value = 0.1868069867954431
$ws.Range("C9").Value = 1.037514270990646
$ws.Range("D9").Value = 3.290947412145117
$ws.Range("E9").Value = 1.814096858534603
$ws.Range("F9").Value = 1.835300565959116
$ws.Range("G9").Value = 30

$ws.Range("B10").Value = 0.2405366504582205
$ws.Range("C10").Value = 1.068640860714179
$ws.Range("D10").Value = 3.366469277036136
$ws.Range("E10").Value = 1.834794069381122
$ws.Range("F10").Value = 1.851155339123688
$ws.Range("G10").Value = 29

$ws.Range("B11").Value = 0.1965726553058998
$ws.Range("C11").Value = 1.065268685810892
$ws.Range("D11").Value = 3.497083030370647
$ws.Range("E11").Value = 1.870048937961423
$ws.Range("F11").Value = 1.893814318830355
$ws.Range("G11").Value = 28
